# Auto-generated Excel COM-interop script to apply the Twintania_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 job sheets
$wb = $excel.ActiveWorkbook

# --- ALC row 53 (Leve Item ID 5479) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 474.9091
$ws.Range("I53").Value = 473.64517
$ws.Range("K53").Value = 473.64517
$ws.Range("M53").Value = 163.35483

# --- ALC row 114 (Leve Item ID 25959) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# --- ALC row 132 (Leve Item ID 44049) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2236.2703
$ws.Range("I132").Value = 2298.6572
$ws.Range("K132").Value = 6895.971600000001
$ws.Range("M132").Value = -4365.971600000001

# --- ALC row 138 (Leve Item ID 44169) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3687.5
$ws.Range("J138").Value = 3563.3704
$ws.Range("L138").Value = 10690.1112
$ws.Range("N138").Value = -20970.1112

# --- ALC row 139 (Leve Item ID 42306) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 97000
$ws.Range("J139").Value = 97000
$ws.Range("L139").Value = 97000
$ws.Range("N139").Value = -107280

# --- ARM row 45 (Leve Item ID 27714) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7415.3335
$ws.Range("J45").Value = 1260.6666
$ws.Range("L45").Value = 1260.6666
$ws.Range("N45").Value = -2014.6666

# --- ARM row 61 (Leve Item ID 43999) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4238.65
$ws.Range("J61").Value = 9398.286
$ws.Range("L61").Value = 9398.286
$ws.Range("N61").Value = -9822.286

# --- ARM row 122 (Leve Item ID 36168) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4793.9375
$ws.Range("I122").Value = 4793.9375
$ws.Range("K122").Value = 14381.8125
$ws.Range("M122").Value = -11931.8125

# --- ARM row 123 (Leve Item ID 34107) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 61683.75
$ws.Range("J123").Value = 61683.75
$ws.Range("L123").Value = 61683.75
$ws.Range("N123").Value = -71483.75

# --- ARM row 132 (Leve Item ID 43997) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11076.0625
$ws.Range("I132").Value = 10495.305
$ws.Range("K132").Value = 31485.915
$ws.Range("M132").Value = -28955.915

# --- ARM row 136 (Leve Item ID 43999) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4238.65
$ws.Range("J136").Value = 9398.286
$ws.Range("L136").Value = 28194.858
$ws.Range("N136").Value = -33294.858

# --- BSM row 64 (Leve Item ID 14184) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1712.3334
$ws.Range("I64").Value = 2481
$ws.Range("J64").Value = 751.5
$ws.Range("K64").Value = 2481
$ws.Range("L64").Value = 751.5
$ws.Range("M64").Value = -2256
$ws.Range("N64").Value = -1201.5

# --- BSM row 67 (Leve Item ID 14184) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1712.3334
$ws.Range("I67").Value = 2481
$ws.Range("J67").Value = 751.5
$ws.Range("K67").Value = 2481
$ws.Range("L67").Value = 751.5
$ws.Range("M67").Value = -1701
$ws.Range("N67").Value = -2311.5

# --- BSM row 95 (Leve Item ID 18194) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 24624
$ws.Range("J95").Value = 24624
$ws.Range("L95").Value = 24624
$ws.Range("N95").Value = -30116

# --- BSM row 134 (Leve Item ID 43998) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10002.186
$ws.Range("I134").Value = 7050.7144
$ws.Range("K134").Value = 21152.1432
$ws.Range("M134").Value = -18617.1432

# --- CRP row 31 (Leve Item ID 44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6550.6665
$ws.Range("I31").Value = 3101.3333
$ws.Range("K31").Value = 3101.3333
$ws.Range("M31").Value = -2806.3333

# --- CRP row 34 (Leve Item ID 44023) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6550.6665
$ws.Range("I34").Value = 3101.3333
$ws.Range("K34").Value = 3101.3333
$ws.Range("M34").Value = -2899.3333

# --- CRP row 41 (Leve Item ID 1917) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20799.8
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30856

# --- CRP row 50 (Leve Item ID 1862) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29066.533
$ws.Range("J50").Value = 29066.533
$ws.Range("L50").Value = 29066.533
$ws.Range("N50").Value = -30316.533

# --- CRP row 59 (Leve Item ID 1942) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 65984.39
$ws.Range("I59").Value = 32534
$ws.Range("J59").Value = 72674.47
$ws.Range("K59").Value = 32534
$ws.Range("L59").Value = 72674.47
$ws.Range("M59").Value = -31389
$ws.Range("N59").Value = -74964.47

# --- CRP row 60 (Leve Item ID 1937) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 9046.333000000001
$ws.Range("I60").Value = 9046.333000000001
$ws.Range("K60").Value = 9046.333000000001
$ws.Range("M60").Value = -8535.333000000001

# --- CUL row 12 (Leve Item ID 4854) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 61.157894
$ws.Range("I12").Value = 93.40000000000001
$ws.Range("J12").Value = 49.642857
$ws.Range("K12").Value = 280.2
$ws.Range("L12").Value = 148.928571
$ws.Range("M12").Value = -107.2
$ws.Range("N12").Value = -494.928571

# --- CUL row 134 (Leve Item ID 44074) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1065.8
$ws.Range("I134").Value = 1065.8
$ws.Range("K134").Value = 3197.4
$ws.Range("M134").Value = 1872.6

# --- CUL row 139 (Leve Item ID 44102) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1007.6667
$ws.Range("I139").Value = 610.7143
$ws.Range("K139").Value = 1832.1429
$ws.Range("M139").Value = 3307.8571

# --- CUL row 140 (Leve Item ID 44097) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2647.125
$ws.Range("I140").Value = 2294.25
$ws.Range("K140").Value = 6882.75
$ws.Range("M140").Value = -1702.75

# --- GSM row 132 (Leve Item ID 44008) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1514.75
$ws.Range("I132").Value = 1514.75
$ws.Range("K132").Value = 4544.25
$ws.Range("M132").Value = -2014.25

# --- LTW row 22 (Leve Item ID 5277) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2385.7273
$ws.Range("I22").Value = 2227.25
$ws.Range("J22").Value = 2476.2856
$ws.Range("K22").Value = 2227.25
$ws.Range("L22").Value = 2476.2856
$ws.Range("M22").Value = -1932.25
$ws.Range("N22").Value = -3066.2856

# --- LTW row 27 (Leve Item ID 5277) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2385.7273
$ws.Range("I27").Value = 2227.25
$ws.Range("J27").Value = 2476.2856
$ws.Range("K27").Value = 2227.25
$ws.Range("L27").Value = 2476.2856
$ws.Range("M27").Value = -2120.25
$ws.Range("N27").Value = -2690.2856

# --- LTW row 40 (Leve Item ID 36248) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1820.6
$ws.Range("I40").Value = 1868
$ws.Range("J40").Value = 1749.5
$ws.Range("K40").Value = 1868
$ws.Range("L40").Value = 1749.5
$ws.Range("M40").Value = -1732
$ws.Range("N40").Value = -2021.5

# --- LTW row 140 (Leve Item ID 42503) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 52367.285
$ws.Range("J140").Value = 52367.285
$ws.Range("L140").Value = 52367.285
$ws.Range("N140").Value = -62727.285

# --- WVR row 81 (Leve Item ID 12596) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1893.5
$ws.Range("J81").Value = 1999
$ws.Range("L81").Value = 3998
$ws.Range("N81").Value = -6120

# --- WVR row 84 (Leve Item ID 12596) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1893.5
$ws.Range("J84").Value = 1999
$ws.Range("L84").Value = 19990
$ws.Range("N84").Value = -30598

# --- WVR row 123 (Leve Item ID 34127) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

# --- WVR row 126 (Leve Item ID 36210) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5391.1284
$ws.Range("I126").Value = 5572.7334
$ws.Range("K126").Value = 16718.2002
$ws.Range("M126").Value = -14248.2002

# --- WVR row 132 (Leve Item ID 44029) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17430.492
$ws.Range("I132").Value = 12870.163
$ws.Range("J132").Value = 31396.5
$ws.Range("K132").Value = 38610.489
$ws.Range("L132").Value = 94189.5
$ws.Range("M132").Value = -36080.489
$ws.Range("N132").Value = -99249.5

# --- WVR row 136 (Leve Item ID 44031) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4765465.5
$ws.Range("I136").Value = 5266672.5
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 15800017.5
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -15797467.5
$ws.Range("N136").Value = -17100
